$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
@("ECs","Cfh","Itgam","FAPs",3,1,0.6718046666666666,2.015414,0.006062562609515538,0.006345114347523019,3,1,0.9636300000000001,2.89089,0.00463428088190967,0.004636691828827675,0.64737113094,5.82634017846,([double]"2.809561799665825E-05"),([double]"2.942033984813722E-05")),
@("ECs","Cfh","Itgam","M1",3,1,0.6718046666666666,2.015414,0.006062562609515538,0.006345114347523019,3,1,71.62512099999999,214.875363,0.3444588990740914,0.344638101013349,48.11809053836465,433.0628148452819,0.002088303642041473,0.002186768159442888),
@("ECs","Cfh","Itgam","M2",3,1,0.6718046666666666,2.015414,0.006062562609515538,0.006345114347523019,3,1,135.0220616666666,405.066185,0.6493469059886695,0.6496847234325412,90.70845113062109,816.3760601755898,0.003936706272851509,0.004122323860018342),
@("ECs","Cfh","Itgam","sCs",3,1,0.6718046666666666,2.015414,0.006062562609515538,0.006345114347523019,1,0.5,0.324361,0.648722,0.001559914055329434,0.001040483725282092,0.2179072334846666,1.307443400908,([double]"9.457076625897975E-06"),([double]"6.601988213651601E-06")),
@("FAPs","Cfh","Itgam","FAPs",3,1,43.16235166666667,129.487055,0.3895097374828606,0.4076632247761514,3,1,0.9636300000000001,2.89089,0.00463428088190967,0.004636691828827675,41.59253693655,374.33283242895,0.001805097529734475,0.001890208743233121),
@("FAPs","Cfh","Itgam","M1",3,1,43.16235166666667,129.487055,0.3895097374828606,0.4076632247761514,3,1,71.62512099999999,214.875363,0.3444588990740914,0.344638101013349,3091.508660769551,27823.57794692597,0.1341700953519845,0.1404962796398309),
@("FAPs","Cfh","Itgam","M2",3,1,43.16235166666667,129.487055,0.3895097374828606,0.4076632247761514,3,1,135.0220616666666,405.066185,0.6493469059886695,0.6496847234325412,5827.869708415018,52450.82737573517,0.2529269428869544,0.2648525694423118),
@("FAPs","Cfh","Itgam","sCs",3,1,43.16235166666667,129.487055,0.3895097374828606,0.4076632247761514,1,0.5,0.324361,0.648722,0.001559914055329434,0.001040483725282092,14.00018354895167,84.00110129371001,0.0006076017141871922,0.0004241669507756008),
@("M1","Cfh","Itgam","FAPs",3,1,21.686232,65.058696,0.1957029264426257,0.2048238552578966,3,1,0.9636300000000001,2.89089,0.00463428088190967,0.004636691828827675,20.89750374216,188.07753367944,0.0009069423305468348,0.0009497050960232717),
@("M1","Cfh","Itgam","M1",3,1,21.686232,65.058696,0.1957029264426257,0.2048238552578966,3,1,71.62512099999999,214.875363,0.3444588990740914,0.344638101013349,1553.278991034072,13979.51091930665,0.06741161458800476,0.07059010451831456),
@("M1","Cfh","Itgam","M2",3,1,21.686232,65.058696,0.1957029264426257,0.2048238552578966,3,1,135.0220616666666,405.066185,0.6493469059886695,0.6496847234325412,2928.11975442164,26353.07778979476,0.1270790897784472,0.1330709297556134),
@("M1","Cfh","Itgam","sCs",3,1,21.686232,65.058696,0.1957029264426257,0.2048238552578966,1,0.5,0.324361,0.648722,0.001559914055329434,0.001040483725282092,7.034167897752001,42.205007386512,0.0003052797456269542,0.0002131158879453763),
@("M2","Cfh","Itgam","FAPs",3,1,30.48803466666666,91.46410399999999,0.2751329786452017,0.2879558237532028,3,1,0.9636300000000001,2.89089,0.00463428088190967,0.004636691828827675,29.37918484584,264.41266361256,0.00127504350291832,0.001335162415059818),
@("M2","Cfh","Itgam","M1",3,1,30.48803466666666,91.46410399999999,0.2751329786452017,0.2879558237532028,3,1,71.62512099999999,214.875363,0.3444588990740914,0.344638101013349,2183.709172052194,19653.38254846975,0.09477200292310169,0.09924054827403844),
@("M2","Cfh","Itgam","M2",3,1,30.48803466666666,91.46410399999999,0.2751329786452017,0.2879558237532028,3,1,135.0220616666666,405.066185,0.6493469059886695,0.6496847234325412,4116.557296858136,37049.01567172323,0.1786567484187084,0.1870804997158892),
@("M2","Cfh","Itgam","sCs",3,1,30.48803466666666,91.46410399999999,0.2751329786452017,0.2879558237532028,1,0.5,0.324361,0.648722,0.001559914055329434,0.001040483725282092,9.889129412514666,59.334776475088,0.0004291838004733031,0.000299613348215406),
@("sCs","Cfh","Itgam","FAPs",2,1,14.8035735,29.607147,0.1335917948197964,0.09321198186522625,3,1,0.9636300000000001,2.89089,0.00463428088190967,0.004636691828827675,14.265167531805,85.59100519083,0.0006191019007133819,0.000432195234663328),
@("sCs","Cfh","Itgam","M1",2,1,14.8035735,29.607147,0.1335917948197964,0.09321198186522625,3,1,71.62512099999999,214.875363,0.3444588990740914,0.344638101013349,1060.307743169893,6361.846459019361,0.04601688256895899,0.0321244004217223),
@("sCs","Cfh","Itgam","M2",2,1,14.8035735,29.607147,0.1335917948197964,0.09321198186522625,3,1,135.0220616666666,405.066185,0.6493469059886695,0.6496847234325412,1998.809014004032,11992.85408402419,0.08674741863170798,0.06055840065870857),
@("sCs","Cfh","Itgam","sCs",2,1,14.8035735,29.607147,0.1335917948197964,0.09321198186522625,1,0.5,0.324361,0.648722,0.001559914055329434,0.001040483725282092,4.8017019040335,19.206807616134,0.0002083917184160863,([double]"9.698555013205742E-05"))
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
